$wb = $excel.ActiveWorkbook

# --- Update sheets LP1912 (1) and 6203-6173 (3): header rows ---
foreach ($idx in 1,3) {
    $ws = $wb.Worksheets.Item($idx)
    $ws.Range("A2").Value = "Última actualización: 05:55:52"
    $ws.Range("A3").Value = "Total filas: 62"

    $ws.Range("A14").Value = "05:55:51"
    $ws.Range("B14").Value = "06:50"
    $ws.Range("C14").Value = "215A_EL PATO"
    $ws.Range("D14").Value = 55
    $ws.Range("A15").Value = "05:55:51"
    $ws.Range("B15").Value = "06:53"
    $ws.Range("C15").Value = "14_ABASTO"
    $ws.Range("D15").Value = 58
    $ws.Range("A16").Value = "04:44:55"
    $ws.Range("B16").Value = "05:43"
    $ws.Range("C16").Value = "14_ABASTO"
    $ws.Range("D16").Value = 59
    $ws.Range("A17").Value = "03:46:12"
    $ws.Range("B17").Value = "04:46"
    $ws.Range("C17").Value = "215A_EL PATO"
    $ws.Range("D17").Value = 60
    $ws.Range("A18").Value = "05:55:51"
    $ws.Range("B18").Value = "07:00"
    $ws.Range("C18").Value = "16_SANTA ANA"
    $ws.Range("D18").Value = 65
    $ws.Range("A19").Value = "04:44:55"
    $ws.Range("B19").Value = "05:51"
    $ws.Range("C19").Value = "17_ROMERO"
    $ws.Range("D19").Value = 67
    $ws.Range("A20").Value = "01:55:38"
    $ws.Range("B20").Value = "03:02"
    $ws.Range("C20").Value = "15_ABASTO"
    $ws.Range("D20").Value = 67
    $ws.Range("A21").Value = "05:55:51"
    $ws.Range("B21").Value = "07:03"
    $ws.Range("C21").Value = "225_GOMEZ"
    $ws.Range("D21").Value = 68
    $ws.Range("A22").Value = "04:01:06"
    $ws.Range("B22").Value = "05:12"
    $ws.Range("C22").Value = "17_ROMERO"
    $ws.Range("D22").Value = 71
    $ws.Range("A23").Value = "05:55:51"
    $ws.Range("B23").Value = "07:06"
    $ws.Range("C23").Value = "215C_EL PATO"
    $ws.Range("D23").Value = 71
    $ws.Range("A24").Value = "00:46:06"
    $ws.Range("B24").Value = "01:58"
    $ws.Range("C24").Value = "14_ABASTO"
    $ws.Range("D24").Value = 72
    $ws.Range("A25").Value = "04:30:03"
    $ws.Range("B25").Value = "05:44"
    $ws.Range("C25").Value = "14_ABASTO"
    $ws.Range("D25").Value = 74
    $ws.Range("A26").Value = "04:44:55"
    $ws.Range("B26").Value = "06:00"
    $ws.Range("C26").Value = "16_SANTA ANA"
    $ws.Range("D26").Value = 76
    $ws.Range("A27").Value = "05:55:51"
    $ws.Range("B27").Value = "07:13"
    $ws.Range("C27").Value = "14X44_ABASTO"
    $ws.Range("D27").Value = 78
    $ws.Range("A28").Value = "04:44:55"
    $ws.Range("B28").Value = "06:03"
    $ws.Range("C28").Value = "10_OLMOS"
    $ws.Range("D28").Value = 79
    $ws.Range("A29").Value = "05:37:13"
    $ws.Range("B29").Value = "07:01"
    $ws.Range("C29").Value = "16_SANTA ANA"
    $ws.Range("D29").Value = 84
    $ws.Range("A30").Value = "04:44:55"
    $ws.Range("B30").Value = "06:10"
    $ws.Range("C30").Value = "215A_EL PATO"
    $ws.Range("D30").Value = 86
    $ws.Range("A31").Value = "05:37:13"
    $ws.Range("B31").Value = "07:04"
    $ws.Range("C31").Value = "225_GOMEZ"
    $ws.Range("D31").Value = 87
    $ws.Range("A32").Value = "03:46:12"
    $ws.Range("B32").Value = "05:16"
    $ws.Range("C32").Value = "17_ROMERO"
    $ws.Range("D32").Value = 90
    $ws.Range("A33").Value = "05:37:13"
    $ws.Range("B33").Value = "07:07"
    $ws.Range("C33").Value = "215C_EL PATO"
    $ws.Range("D33").Value = 90
    $ws.Range("A34").Value = "04:01:06"
    $ws.Range("B34").Value = "05:32"
    $ws.Range("C34").Value = "81_EL PELIGRO"
    $ws.Range("D34").Value = 91
    $ws.Range("A35").Value = "04:30:03"
    $ws.Range("B35").Value = "06:01"
    $ws.Range("C35").Value = "16_SANTA ANA"
    $ws.Range("D35").Value = 91
    $ws.Range("A36").Value = "02:29:13"
    $ws.Range("B36").Value = "04:01"
    $ws.Range("C36").Value = "81_EL PELIGRO"
    $ws.Range("D36").Value = 92
    $ws.Range("A37").Value = "04:58:02"
    $ws.Range("B37").Value = "06:31"
    $ws.Range("C37").Value = "17X38_ROMERO"
    $ws.Range("D37").Value = 93
    $ws.Range("A38").Value = "05:55:51"
    $ws.Range("B38").Value = "07:28"
    $ws.Range("C38").Value = "14_ABASTO"
    $ws.Range("D38").Value = 93
    $ws.Range("A39").Value = "04:58:02"
    $ws.Range("B39").Value = "06:31"
    $ws.Range("C39").Value = "16_SANTA ANA"
    $ws.Range("D39").Value = 93
    $ws.Range("A40").Value = "04:30:03"
    $ws.Range("B40").Value = "06:04"
    $ws.Range("C40").Value = "10_OLMOS"
    $ws.Range("D40").Value = 94
    $ws.Range("A41").Value = "03:46:12"
    $ws.Range("B41").Value = "05:22"
    $ws.Range("C41").Value = "23_HERNANDEZ"
    $ws.Range("D41").Value = 96
    $ws.Range("A42").Value = "01:22:42"
    $ws.Range("B42").Value = "02:58"
    $ws.Range("C42").Value = "215_ALUAR"
    $ws.Range("D42").Value = 96
    $ws.Range("A43").Value = "05:37:13"
    $ws.Range("B43").Value = "07:14"
    $ws.Range("C43").Value = "14X44_ABASTO"
    $ws.Range("D43").Value = 97
    $ws.Range("A44").Value = "04:44:55"
    $ws.Range("B44").Value = "06:23"
    $ws.Range("C44").Value = "11_ETCHEVERRY"
    $ws.Range("D44").Value = 99
    $ws.Range("A45").Value = "05:55:51"
    $ws.Range("B45").Value = "07:35"
    $ws.Range("C45").Value = "17X38_ROMERO"
    $ws.Range("D45").Value = 100
    $ws.Range("A46").Value = "04:30:03"
    $ws.Range("B46").Value = "06:11"
    $ws.Range("C46").Value = "215A_EL PATO"
    $ws.Range("D46").Value = 101
    $ws.Range("A47").Value = "04:58:02"
    $ws.Range("B47").Value = "06:39"
    $ws.Range("C47").Value = "225_C ROCA-H SUR"
    $ws.Range("D47").Value = 101
    $ws.Range("A48").Value = "05:55:51"
    $ws.Range("B48").Value = "07:36"
    $ws.Range("C48").Value = "27_EL RETIRO"
    $ws.Range("D48").Value = 101
    $ws.Range("A49").Value = "04:44:55"
    $ws.Range("B49").Value = "06:26"
    $ws.Range("C49").Value = "23_HERNANDEZ"
    $ws.Range("D49").Value = 102
    $ws.Range("A50").Value = "05:37:13"
    $ws.Range("B50").Value = "07:21"
    $ws.Range("C50").Value = "215A_EL PATO"
    $ws.Range("D50").Value = 104
    $ws.Range("A51").Value = "04:01:06"
    $ws.Range("B51").Value = "05:45"
    $ws.Range("C51").Value = "14_ABASTO"
    $ws.Range("D51").Value = 104
    $ws.Range("A52").Value = "04:44:55"
    $ws.Range("B52").Value = "06:30"
    $ws.Range("C52").Value = "17X38_ROMERO"
    $ws.Range("D52").Value = 106
    $ws.Range("A53").Value = "04:44:55"
    $ws.Range("B53").Value = "06:30"
    $ws.Range("C53").Value = "16_SANTA ANA"
    $ws.Range("D53").Value = 106
    $ws.Range("A54").Value = "05:55:51"
    $ws.Range("B54").Value = "07:43"
    $ws.Range("C54").Value = "10_OLMOS"
    $ws.Range("D54").Value = 108
    $ws.Range("A55").Value = "03:46:12"
    $ws.Range("B55").Value = "05:35"
    $ws.Range("C55").Value = "215B_EL PATO"
    $ws.Range("D55").Value = 109
    $ws.Range("A56").Value = "05:37:13"
    $ws.Range("B56").Value = "07:27"
    $ws.Range("C56").Value = "215A_LA PLATA"
    $ws.Range("D56").Value = 110
    $ws.Range("A57").Value = "04:01:06"
    $ws.Range("B57").Value = "05:52"
    $ws.Range("C57").Value = "17_ROMERO"
    $ws.Range("D57").Value = 111
    $ws.Range("A58").Value = "03:00:53"
    $ws.Range("B58").Value = "04:53"
    $ws.Range("C58").Value = "11_ETCHEVERRY"
    $ws.Range("D58").Value = 113
    $ws.Range("A59").Value = "04:58:02"
    $ws.Range("B59").Value = "06:51"
    $ws.Range("C59").Value = "215A_EL PATO"
    $ws.Range("D59").Value = 113
    $ws.Range("A60").Value = "01:55:38"
    $ws.Range("B60").Value = "03:48"
    $ws.Range("C60").Value = "14_ABASTO"
    $ws.Range("D60").Value = 113
    $ws.Range("A61").Value = "04:44:55"
    $ws.Range("B61").Value = "06:38"
    $ws.Range("C61").Value = "225_C ROCA-H SUR"
    $ws.Range("D61").Value = 114
    $ws.Range("A62").Value = "04:30:03"
    $ws.Range("B62").Value = "06:24"
    $ws.Range("C62").Value = "11_ETCHEVERRY"
    $ws.Range("D62").Value = 114
    $ws.Range("A63").Value = "04:58:02"
    $ws.Range("B63").Value = "06:54"
    $ws.Range("C63").Value = "14_ABASTO"
    $ws.Range("D63").Value = 116
    $ws.Range("A64").Value = "05:37:13"
    $ws.Range("B64").Value = "07:33"
    $ws.Range("C64").Value = "23_HERNANDEZ"
    $ws.Range("D64").Value = 116
    $ws.Range("A65").Value = "05:55:51"
    $ws.Range("B65").Value = "07:51"
    $ws.Range("C65").Value = "15_ABASTO"
    $ws.Range("D65").Value = 116
    $ws.Range("A66").Value = "04:30:03"
    $ws.Range("B66").Value = "06:27"
    $ws.Range("C66").Value = "23_HERNANDEZ"
    $ws.Range("D66").Value = 117
    $ws.Range("A67").Value = "02:47:42"
    $ws.Range("B67").Value = "04:45"
    $ws.Range("C67").Value = "215A_EL PATO"
    $ws.Range("D67").Value = 118
}
# --- Update sheet LP1912-215 (2) ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = "Última actualización: 05:55:52"
$ws2.Range("A3").Value = "Total filas: 15"

$ws2.Range("A9").Value = "05:55:51"
$ws2.Range("B9").Value = "06:50"
$ws2.Range("C9").Value = "215A_EL PATO"
$ws2.Range("D9").Value = 55
$ws2.Range("A10").Value = "03:46:12"
$ws2.Range("B10").Value = "04:46"
$ws2.Range("C10").Value = "215A_EL PATO"
$ws2.Range("D10").Value = 60
$ws2.Range("A11").Value = "05:55:51"
$ws2.Range("B11").Value = "07:06"
$ws2.Range("C11").Value = "215C_EL PATO"
$ws2.Range("D11").Value = 71
$ws2.Range("A12").Value = "04:44:55"
$ws2.Range("B12").Value = "06:10"
$ws2.Range("C12").Value = "215A_EL PATO"
$ws2.Range("D12").Value = 86
$ws2.Range("A13").Value = "05:37:13"
$ws2.Range("B13").Value = "07:07"
$ws2.Range("C13").Value = "215C_EL PATO"
$ws2.Range("D13").Value = 90
$ws2.Range("A14").Value = "01:22:42"
$ws2.Range("B14").Value = "02:58"
$ws2.Range("C14").Value = "215_ALUAR"
$ws2.Range("D14").Value = 96
$ws2.Range("A15").Value = "04:30:03"
$ws2.Range("B15").Value = "06:11"
$ws2.Range("C15").Value = "215A_EL PATO"
$ws2.Range("D15").Value = 101
$ws2.Range("A16").Value = "05:37:13"
$ws2.Range("B16").Value = "07:21"
$ws2.Range("C16").Value = "215A_EL PATO"
$ws2.Range("D16").Value = 104
$ws2.Range("A17").Value = "03:46:12"
$ws2.Range("B17").Value = "05:35"
$ws2.Range("C17").Value = "215B_EL PATO"
$ws2.Range("D17").Value = 109
$ws2.Range("A18").Value = "05:37:13"
$ws2.Range("B18").Value = "07:27"
$ws2.Range("C18").Value = "215A_LA PLATA"
$ws2.Range("D18").Value = 110
$ws2.Range("A19").Value = "04:58:02"
$ws2.Range("B19").Value = "06:51"
$ws2.Range("C19").Value = "215A_EL PATO"
$ws2.Range("D19").Value = 113
$ws2.Range("A20").Value = "02:47:42"
$ws2.Range("B20").Value = "04:45"
$ws2.Range("C20").Value = "215A_EL PATO"
$ws2.Range("D20").Value = 118